# Regenerate s_vals data to filter save games.
# Updates columns B:E and G for rows 2-8 (column F / "Win" is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01253208636536152, 2919.202174992006, 186123.597850132, 2797.565817734744, 191840.3783749451)
    3 = @(1.445647641019636, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 5.507293877332936)
    4 = @(1.445647641019636, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 2.433531715253719)
    5 = @(1.445647641019636, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 5.507293877332936)
    6 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    7 = @(3.272327238179451, 109.9114832445916, 3.223369029078222, 13.86384647080068, 130.27102598265)
    8 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
